$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.370.76'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.849.29'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("D4").Value = '''0.9998'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''240.60'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").Value = '''0.6279'
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''0.07612'
$ws.Range("E8").Value = '  -0.84%  '
$ws.Range("D9").Value = '''0.2914'
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("D10").Value = '''24.64'
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("D11").Value = '''0.07746'
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = '''5.023'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '''0.6799'
$ws.Range("E13").Value = '  -0.23%  '
$ws.Range("E14").Value = '  -5.04%  '
$ws.Range("D15").Value = '''83.11'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = '''6.130'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = '29.391.79'
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").Value = '''228.85'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").Value = '''12.34'
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("D20").Value = '''1.001'
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = '''7.484'
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '''158.73'
$ws.Range("E23").Value = '  +1.09%  '
$ws.Range("D24").Value = '''0.1388'
$ws.Range("D25").Value = '''8.444'
$ws.Range("E25").Value = '  +0.66%  '
$ws.Range("D26").Value = '''17.69'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = '''1.456'
$ws.Range("E27").Value = '  +10.71%  '
$ws.Range("D28").Value = '''1.476'
$ws.Range("E28").Value = '  +0.62%  '
$ws.Range("D29").Value = '''0.05601'
$ws.Range("E29").Value = '  -1.94%  '
$ws.Range("D30").Value = '''4.110'
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("D31").Value = '''4.072'
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").Value = '''1.834'
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("D33").Value = '''1.157'
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("D34").Value = '''0.7009'
$ws.Range("E34").Value = '  -0.97%  '
$ws.Range("D35").Value = '''2.587'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").Value = '1.234.80'
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").Value = '''0.01805'
$ws.Range("E37").Value = '  +0.83%  '
$ws.Range("D38").Value = '''2.732'
$ws.Range("E38").Value = '  -1.57%  '
$ws.Range("D39").Value = '''6.420'
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("D40").Value = '''0.9047'
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("D41").Value = '''1.0000'
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").Value = '''101.50'
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("D43").Value = '''65.56'
$ws.Range("E43").Value = '  -1.12%  '
$ws.Range("D44").Value = '''7.202'
$ws.Range("E44").Value = '  +1.11%  '
$ws.Range("D45").Value = '''0.3997'
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("D46").Value = '''9.047'
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '''0.1152'
$ws.Range("E47").Value = '  +1.78%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''1.685'
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("E49").Value = '  -3.36%  '
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("D51").Value = '''0.4630'
$ws.Range("E51").Value = '  +0.06%  '
